$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "28.735.04", "  +6.88%  ")
    ,@("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.812.63", "  +4.96%  ")
    ,@("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "'0.9990", "  +0.18%  ")
    ,@("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "'250.98", "  +3.70%  ")
    ,@("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "'0.9993", "  +0.16%  ")
    ,@("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "'0.4968", "  +1.54%  ")
    ,@("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "'0.2788", "  +7.54%  ")
    ,@("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "'0.06381", "  +2.67%  ")
    ,@("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.805.55", "  +4.46%  ")
    ,@("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "'16.74", "  +4.66%  ")
    ,@("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "'0.07114", "  +3.07%  ")
    ,@("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "'0.6484", "  +6.44%  ")
    ,@("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "'4.706", "  +4.95%  ")
    ,@("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "'81.86", "  +5.95%  ")
    ,@("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "28.698.64", "  +6.83%  ")
    ,@("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "'0.9990", "  +0.08%  ")
    ,@("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "'0.000007384", "  +2.97%  ")
    ,@("BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "'0.9989", "  +0.16%  ")
    ,@("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "'12.30", "  +7.54%  ")
    ,@("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.038.88", "  +4.36%  ")
    ,@("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "'4.622", "  +4.42%  ")
    ,@("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "'8.898", "  +3.78%  ")
    ,@("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "'5.322", "  +4.20%  ")
    ,@("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "'142.66", "  +3.02%  ")
    ,@("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "'16.03", "  +4.79%  ")
    ,@("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "'1.881", "  +5.74%  ")
    ,@("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "'112.74", "  +6.27%  ")
    ,@("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "'1.390", "  +0.76%  ")
    ,@("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "'4.171", "  +5.67%  ")
    ,@("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "'0.08362", "  +4.55%  ")
    ,@("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "'3.838", "  +4.04%  ")
    ,@("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "'0.04972", "  +9.71%  ")
    ,@("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "'1.088", "  +8.02%  ")
    ,@("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "'0.6769", "  +8.23%  ")
    ,@("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "'2.663", "  +2.61%  ")
    ,@("MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "'2.683", "  +9.46%  ")
    ,@("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "'0.9579", "  +2.28%  ")
    ,@("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "'2.136", "  +3.63%  ")
    ,@("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "'0.01590", "  +5.83%  ")
    ,@("FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "'5.959", "  +5.47%  ")
    ,@("PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "'0.9997", "  +0.21%  ")
    ,@("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "'101.02", "  +1.67%  ")
    ,@("TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "'0.4112", "  +6.70%  ")
    ,@("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "'7.208", "  +4.52%  ")
    ,@("Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "'0.1224", "  +5.26%  ")
    ,@("Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "'0.05498", "  +1.96%  ")
    ,@("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "'8.216", "  +3.97%  ")
    ,@("Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "'31.42", "  +4.14%  ")
    ,@("Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "'0.3621", "  +7.14%  ")
    ,@("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "'1.303", "  +5.34%  ")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}
